# Auto-generated edit script: updates numeric leve-profit columns (H..N)
# across multiple worksheets, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 731.39026  # H15: 1090.4884 -> 731.39026
$ws.Cells.Item(15, 9).Value = 731.39026  # I15: 1090.4884 -> 731.39026
$ws.Cells.Item(15, 11).Value = 2194.17078  # K15: 3271.4652 -> 2194.17078
$ws.Cells.Item(15, 13).Value = -2025.17078  # M15: -3102.4652 -> -2025.17078

$ws.Cells.Item(33, 8).Value = 118.8  # H33: 112.57143 -> 118.8
$ws.Cells.Item(33, 9).Value = 123.5  # I33: 117.8 -> 123.5
$ws.Cells.Item(33, 10).Value = 100  # J33: 99.5 -> 100
$ws.Cells.Item(33, 11).Value = 123.5  # K33: 117.8 -> 123.5
$ws.Cells.Item(33, 12).Value = 100  # L33: 99.5 -> 100
$ws.Cells.Item(33, 13).Value = 105.5  # M33: 111.2 -> 105.5
$ws.Cells.Item(33, 14).Value = -558  # N33: -557.5 -> -558

$ws.Cells.Item(62, 8).Value = 2569.2856  # H62: 2698.4 -> 2569.2856
$ws.Cells.Item(62, 9).Value = 2397.8  # I62: 2399.6 -> 2397.8
$ws.Cells.Item(62, 10).Value = 2998  # J62: 2997.2 -> 2998
$ws.Cells.Item(62, 11).Value = 2397.8  # K62: 2399.6 -> 2397.8
$ws.Cells.Item(62, 12).Value = 2998  # L62: 2997.2 -> 2998
$ws.Cells.Item(62, 13).Value = -1773.8  # M62: -1775.6 -> -1773.8
$ws.Cells.Item(62, 14).Value = -4246  # N62: -4245.2 -> -4246

$ws.Cells.Item(65, 8).Value = 2569.2856  # H65: 2698.4 -> 2569.2856
$ws.Cells.Item(65, 9).Value = 2397.8  # I65: 2399.6 -> 2397.8
$ws.Cells.Item(65, 10).Value = 2998  # J65: 2997.2 -> 2998
$ws.Cells.Item(65, 11).Value = 11989  # K65: 11998 -> 11989
$ws.Cells.Item(65, 12).Value = 14990  # L65: 14986 -> 14990
$ws.Cells.Item(65, 13).Value = -8869  # M65: -8878 -> -8869
$ws.Cells.Item(65, 14).Value = -21230  # N65: -21226 -> -21230

$ws.Cells.Item(98, 8).Value = 2652.889  # H98: 1865.9584 -> 2652.889
$ws.Cells.Item(98, 9).Value = 2859.5625  # I98: 1965.9048 -> 2859.5625
$ws.Cells.Item(98, 10).Value = 999.5  # J98: 1166.3334 -> 999.5
$ws.Cells.Item(98, 11).Value = 2859.5625  # K98: 1965.9048 -> 2859.5625
$ws.Cells.Item(98, 12).Value = 999.5  # L98: 1166.3334 -> 999.5
$ws.Cells.Item(98, 13).Value = -1361.5625  # M98: -467.9048 -> -1361.5625
$ws.Cells.Item(98, 14).Value = -3995.5  # N98: -4162.3334 -> -3995.5

$ws.Cells.Item(113, 8).Value = 42950.6  # H113: 36025.5 -> 42950.6
$ws.Cells.Item(113, 9).Value = 42950.6  # I113: 36025.5 -> 42950.6
$ws.Cells.Item(113, 11).Value = 42950.6  # K113: 36025.5 -> 42950.6
$ws.Cells.Item(113, 13).Value = -39696.6  # M113: -32771.5 -> -39696.6

$ws.Cells.Item(122, 8).Value = 2652.889  # H122: 1865.9584 -> 2652.889
$ws.Cells.Item(122, 9).Value = 2859.5625  # I122: 1965.9048 -> 2859.5625
$ws.Cells.Item(122, 10).Value = 999.5  # J122: 1166.3334 -> 999.5
$ws.Cells.Item(122, 11).Value = 8578.6875  # K122: 5897.7144 -> 8578.6875
$ws.Cells.Item(122, 12).Value = 2998.5  # L122: 3499.0002 -> 2998.5
$ws.Cells.Item(122, 13).Value = -6128.6875  # M122: -3447.7144 -> -6128.6875
$ws.Cells.Item(122, 14).Value = -7898.5  # N122: -8399.0002 -> -7898.5

$ws.Cells.Item(137, 8).Value = 1433.88  # H137: 1292.0968 -> 1433.88
$ws.Cells.Item(137, 9).Value = 1365.7727  # I137: 1217.5 -> 1365.7727
$ws.Cells.Item(137, 10).Value = 1933.3334  # J137: 1680 -> 1933.3334
$ws.Cells.Item(137, 11).Value = 4097.3181  # K137: 3652.5 -> 4097.3181
$ws.Cells.Item(137, 12).Value = 5800.0002  # L137: 5040 -> 5800.0002
$ws.Cells.Item(137, 13).Value = -1547.3181  # M137: -1102.5 -> -1547.3181
$ws.Cells.Item(137, 14).Value = -10900.0002  # N137: -10140 -> -10900.0002

$ws.Cells.Item(138, 8).Value = 2726.762  # H138: 2742.8462 -> 2726.762
$ws.Cells.Item(138, 10).Value = 3597.5833  # J138: 3957.5557 -> 3597.5833
$ws.Cells.Item(138, 12).Value = 10792.7499  # L138: 11872.6671 -> 10792.7499
$ws.Cells.Item(138, 14).Value = -21072.7499  # N138: -22152.6671 -> -21072.7499

$ws.Cells.Item(139, 8).Value = 48550  # H139: 48500 -> 48550
$ws.Cells.Item(139, 10).Value = 48550  # J139: 48500 -> 48550
$ws.Cells.Item(139, 12).Value = 48550  # L139: 48500 -> 48550
$ws.Cells.Item(139, 14).Value = -58830  # N139: -58780 -> -58830

$ws.Cells.Item(141, 8).Value = 1002195.44  # H141: 1079257.5 -> 1002195.44
$ws.Cells.Item(141, 9).Value = 1274176.9  # I141: 1401555.8 -> 1274176.9
$ws.Cells.Item(141, 11).Value = 3822530.7  # K141: 4204667.4 -> 3822530.7
$ws.Cells.Item(141, 13).Value = -3817350.7  # M141: -4199487.4 -> -3817350.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2542.4482  # H32: 2599.305 -> 2542.4482
$ws.Cells.Item(32, 9).Value = 1829.8701  # I32: 1845.1389 -> 1829.8701
$ws.Cells.Item(32, 11).Value = 1829.8701  # K32: 1845.1389 -> 1829.8701
$ws.Cells.Item(32, 13).Value = -1542.8701  # M32: -1558.1389 -> -1542.8701

$ws.Cells.Item(45, 8).Value = 1812.5883  # H45: 1847.25 -> 1812.5883
$ws.Cells.Item(45, 9).Value = 1758  # I45: 1828.8334 -> 1758
$ws.Cells.Item(45, 10).Value = 1874  # J45: 1865.6666 -> 1874
$ws.Cells.Item(45, 11).Value = 1758  # K45: 1828.8334 -> 1758
$ws.Cells.Item(45, 12).Value = 1874  # L45: 1865.6666 -> 1874
$ws.Cells.Item(45, 13).Value = -1381  # M45: -1451.8334 -> -1381
$ws.Cells.Item(45, 14).Value = -2628  # N45: -2619.6666 -> -2628

$ws.Cells.Item(61, 8).Value = 7953  # H61: 6301.8335 -> 7953
$ws.Cells.Item(61, 10).Value = 12500  # J61: 7749.75 -> 12500
$ws.Cells.Item(61, 12).Value = 12500  # L61: 7749.75 -> 12500
$ws.Cells.Item(61, 14).Value = -12924  # N61: -8173.75 -> -12924

$ws.Cells.Item(122, 8).Value = 1300  # H122: 1275 -> 1300
$ws.Cells.Item(122, 9).Value = 1300  # I122: 1275 -> 1300
$ws.Cells.Item(122, 11).Value = 3900  # K122: 3825 -> 3900
$ws.Cells.Item(122, 13).Value = -1450  # M122: -1375 -> -1450

$ws.Cells.Item(132, 8).Value = 3399.7  # H132: 3023.6155 -> 3399.7
$ws.Cells.Item(132, 9).Value = 4699.5  # I132: 3177.25 -> 4699.5
$ws.Cells.Item(132, 10).Value = 3074.75  # J132: 2955.3333 -> 3074.75
$ws.Cells.Item(132, 11).Value = 14098.5  # K132: 9531.75 -> 14098.5
$ws.Cells.Item(132, 12).Value = 9224.25  # L132: 8865.999899999999 -> 9224.25
$ws.Cells.Item(132, 13).Value = -11568.5  # M132: -7001.75 -> -11568.5
$ws.Cells.Item(132, 14).Value = -14284.25  # N132: -13925.9999 -> -14284.25

$ws.Cells.Item(136, 8).Value = 7953  # H136: 6301.8335 -> 7953
$ws.Cells.Item(136, 10).Value = 12500  # J136: 7749.75 -> 12500
$ws.Cells.Item(136, 12).Value = 37500  # L136: 23249.25 -> 37500
$ws.Cells.Item(136, 14).Value = -42600  # N136: -28349.25 -> -42600

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 334.375  # H94: 298 -> 334.375
$ws.Cells.Item(94, 9).Value = 346.66666  # I94: 308.2069 -> 346.66666
$ws.Cells.Item(94, 11).Value = 346.66666  # K94: 308.2069 -> 346.66666
$ws.Cells.Item(94, 13).Value = 104.33334  # M94: 142.7931 -> 104.33334

$ws.Cells.Item(107, 8).Value = 3049.25  # H107: 3178.25 -> 3049.25
$ws.Cells.Item(107, 9).Value = 3914  # I107: 3566.8333 -> 3914
$ws.Cells.Item(107, 10).Value = 1608  # J107: 2012.5 -> 1608
$ws.Cells.Item(107, 11).Value = 3914  # K107: 3566.8333 -> 3914
$ws.Cells.Item(107, 12).Value = 1608  # L107: 2012.5 -> 1608
$ws.Cells.Item(107, 13).Value = -1994  # M107: -1646.8333 -> -1994
$ws.Cells.Item(107, 14).Value = -5448  # N107: -5852.5 -> -5448

$ws.Cells.Item(134, 8).Value = 9798.883  # H134: 9310.056 -> 9798.883
$ws.Cells.Item(134, 9).Value = 9798.883  # I134: 9310.056 -> 9798.883
$ws.Cells.Item(134, 11).Value = 29396.649  # K134: 27930.168 -> 29396.649
$ws.Cells.Item(134, 13).Value = -26861.649  # M134: -25395.168 -> -26861.649

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1801.4839  # H31: 1919.5714 -> 1801.4839
$ws.Cells.Item(31, 9).Value = 1123.15  # I31: 1197.9412 -> 1123.15
$ws.Cells.Item(31, 11).Value = 1123.15  # K31: 1197.9412 -> 1123.15
$ws.Cells.Item(31, 13).Value = -828.1500000000001  # M31: -902.9412 -> -828.1500000000001

$ws.Cells.Item(34, 8).Value = 1801.4839  # H34: 1919.5714 -> 1801.4839
$ws.Cells.Item(34, 9).Value = 1123.15  # I34: 1197.9412 -> 1123.15
$ws.Cells.Item(34, 11).Value = 1123.15  # K34: 1197.9412 -> 1123.15
$ws.Cells.Item(34, 13).Value = -921.1500000000001  # M34: -995.9412 -> -921.1500000000001

$ws.Cells.Item(86, 8).Value = 1953.909  # H86: 1980.091 -> 1953.909
$ws.Cells.Item(86, 10).Value = 2284.5  # J86: 2332.5 -> 2284.5
$ws.Cells.Item(86, 12).Value = 2284.5  # L86: 2332.5 -> 2284.5
$ws.Cells.Item(86, 14).Value = -4530.5  # N86: -4578.5 -> -4530.5

$ws.Cells.Item(89, 8).Value = 1953.909  # H89: 1980.091 -> 1953.909
$ws.Cells.Item(89, 10).Value = 2284.5  # J89: 2332.5 -> 2284.5
$ws.Cells.Item(89, 12).Value = 11422.5  # L89: 11662.5 -> 11422.5
$ws.Cells.Item(89, 14).Value = -22654.5  # N89: -22894.5 -> -22654.5

$ws.Cells.Item(99, 8).Value = 3304.6667  # H99: 1252164 -> 3304.6667
$ws.Cells.Item(99, 9).Value = 0  # I99: 3334132.8 -> 0
$ws.Cells.Item(99, 10).Value = 3304.6667  # J99: 2982.8 -> 3304.6667
$ws.Cells.Item(99, 11).Value = 0  # K99: 3334132.8 -> 0
$ws.Cells.Item(99, 12).Value = 3304.6667  # L99: 2982.8 -> 3304.6667
$ws.Cells.Item(99, 13).ClearContents()  # M99: remove (was -3332634.8)
$ws.Cells.Item(99, 14).Value = -6300.6667  # N99: -5978.8 -> -6300.6667

$ws.Cells.Item(122, 8).Value = 1646.5294  # H122: 1637.8334 -> 1646.5294
$ws.Cells.Item(122, 9).Value = 1705.6875  # I122: 1693 -> 1705.6875
$ws.Cells.Item(122, 11).Value = 5117.0625  # K122: 5079 -> 5117.0625
$ws.Cells.Item(122, 13).Value = -2667.0625  # M122: -2629 -> -2667.0625

$ws.Cells.Item(126, 8).Value = 3304.6667  # H126: 1252164 -> 3304.6667
$ws.Cells.Item(126, 9).Value = 0  # I126: 3334132.8 -> 0
$ws.Cells.Item(126, 10).Value = 3304.6667  # J126: 2982.8 -> 3304.6667
$ws.Cells.Item(126, 11).Value = 0  # K126: 10002398.4 -> 0
$ws.Cells.Item(126, 12).Value = 9914.000100000001  # L126: 8948.400000000001 -> 9914.000100000001
$ws.Cells.Item(126, 13).ClearContents()  # M126: remove (was -9999928.399999999)
$ws.Cells.Item(126, 14).Value = -14854.0001  # N126: -13888.4 -> -14854.0001

$ws.Cells.Item(132, 8).Value = 2175.5  # H132: 2248.348 -> 2175.5
$ws.Cells.Item(132, 9).Value = 1128.6428  # I132: 1177 -> 1128.6428
$ws.Cells.Item(132, 11).Value = 3385.9284  # K132: 3531 -> 3385.9284
$ws.Cells.Item(132, 13).Value = -855.9284000000002  # M132: -1001 -> -855.9284000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 0  # H37: 100000 -> 0
$ws.Cells.Item(37, 10).Value = 0  # J37: 100000 -> 0
$ws.Cells.Item(37, 12).Value = 0  # L37: 300000 -> 0
$ws.Cells.Item(37, 14).ClearContents()  # N37: remove (was -300224)

$ws.Cells.Item(129, 8).Value = 28780.54  # H129: 28759.424 -> 28780.54
$ws.Cells.Item(129, 9).Value = 469.2857  # I129: 460.625 -> 469.2857
$ws.Cells.Item(129, 10).Value = 39211  # J129: 41336.668 -> 39211
$ws.Cells.Item(129, 11).Value = 1407.8571  # K129: 1381.875 -> 1407.8571
$ws.Cells.Item(129, 12).Value = 117633  # L129: 124010.004 -> 117633
$ws.Cells.Item(129, 13).Value = 3592.1429  # M129: 3618.125 -> 3592.1429
$ws.Cells.Item(129, 14).Value = -127633  # N129: -134010.004 -> -127633

$ws.Cells.Item(140, 8).Value = 1605.1  # H140: 1613.4814 -> 1605.1
$ws.Cells.Item(140, 9).Value = 935.94116  # I140: 932 -> 935.94116
$ws.Cells.Item(140, 10).Value = 2480.1538  # J140: 2604.7273 -> 2480.1538
$ws.Cells.Item(140, 11).Value = 2807.82348  # K140: 2796 -> 2807.82348
$ws.Cells.Item(140, 12).Value = 7440.4614  # L140: 7814.1819 -> 7440.4614
$ws.Cells.Item(140, 13).Value = 2372.17652  # M140: 2384 -> 2372.17652
$ws.Cells.Item(140, 14).Value = -17800.4614  # N140: -18174.1819 -> -17800.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 25000  # H49: 0 -> 25000
$ws.Cells.Item(49, 10).Value = 25000  # J49: 0 -> 25000
$ws.Cells.Item(49, 12).Value = 25000  # L49: 0 -> 25000
$ws.Cells.Item(49, 14).Value = -25368  # N49: new cell

$ws.Cells.Item(102, 8).Value = 5199.5  # H102: 3156.3635 -> 5199.5
$ws.Cells.Item(102, 9).Value = 6133  # I102: 3703.5 -> 6133
$ws.Cells.Item(102, 10).Value = 2399  # J102: 2499.8 -> 2399
$ws.Cells.Item(102, 11).Value = 6133  # K102: 3703.5 -> 6133
$ws.Cells.Item(102, 12).Value = 2399  # L102: 2499.8 -> 2399
$ws.Cells.Item(102, 13).Value = -4511  # M102: -2081.5 -> -4511
$ws.Cells.Item(102, 14).Value = -5643  # N102: -5743.8 -> -5643

$ws.Cells.Item(122, 8).Value = 2752.1667  # H122: 3002.1667 -> 2752.1667
$ws.Cells.Item(122, 9).Value = 2004.3334  # I122: 2006.5 -> 2004.3334
$ws.Cells.Item(122, 11).Value = 6013.0002  # K122: 6019.5 -> 6013.0002
$ws.Cells.Item(122, 13).Value = -3563.0002  # M122: -3569.5 -> -3563.0002

$ws.Cells.Item(126, 8).Value = 1769800  # H126: 1826847.6 -> 1769800
$ws.Cells.Item(126, 9).Value = 2139261.5  # I126: 2224779 -> 2139261.5
$ws.Cells.Item(126, 11).Value = 6417784.5  # K126: 6674337 -> 6417784.5
$ws.Cells.Item(126, 13).Value = -6415314.5  # M126: -6671867 -> -6415314.5

$ws.Cells.Item(132, 8).Value = 7695106.5  # H132: 6413389 -> 7695106.5
$ws.Cells.Item(132, 9).Value = 12822511  # I132: 12822112 -> 12822511
$ws.Cells.Item(132, 10).Value = 3999  # J132: 4665.6665 -> 3999
$ws.Cells.Item(132, 11).Value = 38467533  # K132: 38466336 -> 38467533
$ws.Cells.Item(132, 12).Value = 11997  # L132: 13996.9995 -> 11997
$ws.Cells.Item(132, 13).Value = -38465003  # M132: -38463806 -> -38465003
$ws.Cells.Item(132, 14).Value = -17057  # N132: -19056.9995 -> -17057

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5314.7144  # H7: 6720.8 -> 5314.7144
$ws.Cells.Item(7, 9).Value = 3440.6  # I7: 4534.6665 -> 3440.6
$ws.Cells.Item(7, 11).Value = 3440.6  # K7: 4534.6665 -> 3440.6
$ws.Cells.Item(7, 13).Value = -3328.6  # M7: -4422.6665 -> -3328.6

$ws.Cells.Item(40, 8).Value = 6449  # H40: 7359.8 -> 6449
$ws.Cells.Item(40, 9).Value = 2173.5  # I40: 2266.3333 -> 2173.5
$ws.Cells.Item(40, 11).Value = 2173.5  # K40: 2266.3333 -> 2173.5
$ws.Cells.Item(40, 13).Value = -2037.5  # M40: -2130.3333 -> -2037.5

$ws.Cells.Item(63, 8).Value = 0  # H63: 46985 -> 0
$ws.Cells.Item(63, 10).Value = 0  # J63: 46985 -> 0
$ws.Cells.Item(63, 12).Value = 0  # L63: 46985 -> 0
$ws.Cells.Item(63, 14).ClearContents()  # N63: remove (was -48483)

$ws.Cells.Item(66, 8).Value = 0  # H66: 46985 -> 0
$ws.Cells.Item(66, 10).Value = 0  # J66: 46985 -> 0
$ws.Cells.Item(66, 12).Value = 0  # L66: 140955 -> 0
$ws.Cells.Item(66, 14).ClearContents()  # N66: remove (was -148443)

$ws.Cells.Item(122, 8).Value = 14668  # H122: 11071.857 -> 14668
$ws.Cells.Item(122, 9).Value = 14502  # I122: 10750.75 -> 14502
$ws.Cells.Item(122, 10).Value = 15000  # J122: 11500 -> 15000
$ws.Cells.Item(122, 11).Value = 43506  # K122: 32252.25 -> 43506
$ws.Cells.Item(122, 12).Value = 45000  # L122: 34500 -> 45000
$ws.Cells.Item(122, 13).Value = -41056  # M122: -29802.25 -> -41056
$ws.Cells.Item(122, 14).Value = -49900  # N122: -39400 -> -49900

$ws.Cells.Item(126, 8).Value = 5314.7144  # H126: 6720.8 -> 5314.7144
$ws.Cells.Item(126, 9).Value = 3440.6  # I126: 4534.6665 -> 3440.6
$ws.Cells.Item(126, 11).Value = 10321.8  # K126: 13603.9995 -> 10321.8
$ws.Cells.Item(126, 13).Value = -7851.799999999999  # M126: -11133.9995 -> -7851.799999999999

$ws.Cells.Item(132, 8).Value = 2206.818  # H132: 1738.4474 -> 2206.818
$ws.Cells.Item(132, 9).Value = 2556.1667  # I132: 1502.0625 -> 2556.1667
$ws.Cells.Item(132, 10).Value = 2075.8125  # J132: 1910.3636 -> 2075.8125
$ws.Cells.Item(132, 11).Value = 7668.500100000001  # K132: 4506.1875 -> 7668.500100000001
$ws.Cells.Item(132, 12).Value = 6227.4375  # L132: 5731.0908 -> 6227.4375
$ws.Cells.Item(132, 13).Value = -5138.500100000001  # M132: -1976.1875 -> -5138.500100000001
$ws.Cells.Item(132, 14).Value = -11287.4375  # N132: -10791.0908 -> -11287.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1286.5385  # H81: 1395.5 -> 1286.5385
$ws.Cells.Item(81, 9).Value = 1158.6364  # I81: 1276.6 -> 1158.6364
$ws.Cells.Item(81, 11).Value = 2317.2728  # K81: 2553.2 -> 2317.2728
$ws.Cells.Item(81, 13).Value = -1256.2728  # M81: -1492.2 -> -1256.2728

$ws.Cells.Item(84, 8).Value = 1286.5385  # H84: 1395.5 -> 1286.5385
$ws.Cells.Item(84, 9).Value = 1158.6364  # I84: 1276.6 -> 1158.6364
$ws.Cells.Item(84, 11).Value = 11586.364  # K84: 12766 -> 11586.364
$ws.Cells.Item(84, 13).Value = -6282.364000000001  # M84: -7462 -> -6282.364000000001

$ws.Cells.Item(107, 8).Value = 524.64514  # H107: 524.6129 -> 524.64514
$ws.Cells.Item(107, 9).Value = 417.07407  # I107: 417.03705 -> 417.07407
$ws.Cells.Item(107, 11).Value = 1251.22221  # K107: 1251.11115 -> 1251.22221
$ws.Cells.Item(107, 13).Value = 668.7777900000001  # M107: 668.8888499999998 -> 668.7777900000001

$ws.Cells.Item(113, 8).Value = 577.0714  # H113: 517.8570999999999 -> 577.0714
$ws.Cells.Item(113, 9).Value = 348.25  # I113: 411.53845 -> 348.25
$ws.Cells.Item(113, 10).Value = 1950  # J113: 1900 -> 1950
$ws.Cells.Item(113, 11).Value = 1044.75  # K113: 1234.61535 -> 1044.75
$ws.Cells.Item(113, 12).Value = 5850  # L113: 5700 -> 5850
$ws.Cells.Item(113, 13).Value = 1125.25  # M113: 935.38465 -> 1125.25
$ws.Cells.Item(113, 14).Value = -10190  # N113: -10040 -> -10190

$ws.Cells.Item(126, 8).Value = 4167.8184  # H126: 12209.583 -> 4167.8184
$ws.Cells.Item(126, 9).Value = 1980.75  # I126: 12946.111 -> 1980.75
$ws.Cells.Item(126, 11).Value = 5942.25  # K126: 38838.333 -> 5942.25
$ws.Cells.Item(126, 13).Value = -3472.25  # M126: -36368.333 -> -3472.25

$ws.Cells.Item(132, 8).Value = 1426.1177  # H132: 1026.3214 -> 1426.1177
$ws.Cells.Item(132, 9).Value = 1049.6666  # I132: 778.38464 -> 1049.6666
$ws.Cells.Item(132, 11).Value = 3148.9998  # K132: 2335.15392 -> 3148.9998
$ws.Cells.Item(132, 13).Value = -618.9998000000001  # M132: 194.8460800000003 -> -618.9998000000001

$ws.Cells.Item(136, 8).Value = 2929.087  # H136: 2939.4783 -> 2929.087
$ws.Cells.Item(136, 10).Value = 2218.4666  # J136: 2234.4 -> 2218.4666
$ws.Cells.Item(136, 12).Value = 6655.399800000001  # L136: 6703.200000000001 -> 6655.399800000001
$ws.Cells.Item(136, 14).Value = -11755.3998  # N136: -11803.2 -> -11755.3998
